# Open APC Sweden - Lnu 2019: convert the "is_hybrid" column (E) from text
# TRUE/FALSE (shared-string lookups) to native boolean cells, and update
# the last selected cell on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")

# Mapping row -> is_hybrid boolean value, taken from the original TRUE/FALSE
# text values in column E (rows 2-61).
$values = @{
    2  = $true
    3  = $true
    4  = $false
    5  = $true
    6  = $false
    7  = $true
    8  = $true
    9  = $true
    10 = $true
    11 = $false
    12 = $true
    13 = $false
    14 = $false
    15 = $false
    16 = $false
    17 = $false
    18 = $false
    19 = $false
    20 = $true
    21 = $false
    22 = $true
    23 = $false
    24 = $false
    25 = $false
    26 = $true
    27 = $false
    28 = $false
    29 = $false
    30 = $false
    31 = $false
    32 = $false
    33 = $false
    34 = $false
    35 = $false
    36 = $false
    37 = $false
    38 = $false
    39 = $true
    40 = $true
    41 = $true
    42 = $false
    43 = $false
    44 = $false
    45 = $false
    46 = $false
    47 = $false
    48 = $false
    49 = $false
    50 = $false
    51 = $false
    52 = $false
    53 = $false
    54 = $false
    55 = $false
    56 = $false
    57 = $false
    58 = $true
    59 = $false
    60 = $false
    61 = $true
}

foreach ($row in $values.Keys) {
    $ws.Range("E$row").Value = $values[$row]
}

# Restore the active selection to the cell it was on when the file was
# last saved (D36).
$ws.Range("D36").Select()
